# Apply Asura_Profits market-data refresh (scheduled runner update)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1440.2
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1440.2
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1440.2
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2092.2
$ws.Range("H43").Value = 2625.5386
$ws.Range("I43").Value = 2236.889
$ws.Range("J43").Value = 3500
$ws.Range("K43").Value = 2236.889
$ws.Range("L43").Value = 3500
$ws.Range("M43").Value = -2167.889
$ws.Range("N43").Value = -3638
$ws.Range("H62").Value = 4499.3335
$ws.Range("I62").Value = 4499.3335
$ws.Range("K62").Value = 4499.3335
$ws.Range("M62").Value = -3875.3335
$ws.Range("H65").Value = 4499.3335
$ws.Range("I65").Value = 4499.3335
$ws.Range("K65").Value = 22496.6675
$ws.Range("M65").Value = -19376.6675
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488
$ws.Range("H86").Value = 1792
$ws.Range("I86").Value = 818.1667
$ws.Range("J86").Value = 2522.375
$ws.Range("K86").Value = 818.1667
$ws.Range("L86").Value = 2522.375
$ws.Range("M86").Value = 304.8333
$ws.Range("N86").Value = -4768.375
$ws.Range("H89").Value = 1792
$ws.Range("I89").Value = 818.1667
$ws.Range("J89").Value = 2522.375
$ws.Range("K89").Value = 4090.8335
$ws.Range("L89").Value = 12611.875
$ws.Range("M89").Value = 1525.1665
$ws.Range("N89").Value = -23843.875
$ws.Range("H106").Value = 3476
$ws.Range("I106").Value = 3845
$ws.Range("K106").Value = 3845
$ws.Range("M106").Value = -3214
$ws.Range("H129").Value = 1049.7444
$ws.Range("I129").Value = 309.57144
$ws.Range("J129").Value = 1112.1687
$ws.Range("K129").Value = 928.71432
$ws.Range("L129").Value = 3336.5061
$ws.Range("M129").Value = 4071.28568
$ws.Range("N129").Value = -13336.5061
$ws.Range("H135").Value = 906
$ws.Range("I135").Value = 679.6111
$ws.Range("J135").Value = 1924.75
$ws.Range("K135").Value = 6116.4999
$ws.Range("L135").Value = 17322.75
$ws.Range("M135").Value = -3581.4999
$ws.Range("N135").Value = -22392.75
$ws.Range("H138").Value = 3865.3157
$ws.Range("I138").Value = 2291.9546
$ws.Range("K138").Value = 6875.8638
$ws.Range("M138").Value = -1735.8638
$ws.Range("H141").Value = 4985.4136
$ws.Range("I141").Value = 1595.6296
$ws.Range("J141").Value = 50747.5
$ws.Range("K141").Value = 4786.8888
$ws.Range("L141").Value = 152242.5
$ws.Range("M141").Value = 393.1112000000003
$ws.Range("N141").Value = -162602.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2710.742
$ws.Range("I61").Value = 2597.1304
$ws.Range("J61").Value = 3037.375
$ws.Range("K61").Value = 2597.1304
$ws.Range("L61").Value = 3037.375
$ws.Range("M61").Value = -2385.1304
$ws.Range("N61").Value = -3461.375
$ws.Range("H122").Value = 5141.5
$ws.Range("I122").Value = 6004.905
$ws.Range("J122").Value = 2551.2856
$ws.Range("K122").Value = 18014.715
$ws.Range("L122").Value = 7653.8568
$ws.Range("M122").Value = -15564.715
$ws.Range("N122").Value = -12553.8568
$ws.Range("H132").Value = 5160.6284
$ws.Range("I132").Value = 5244.4375
$ws.Range("J132").Value = 4266.6665
$ws.Range("K132").Value = 15733.3125
$ws.Range("L132").Value = 12799.9995
$ws.Range("M132").Value = -13203.3125
$ws.Range("N132").Value = -17859.9995
$ws.Range("H136").Value = 2710.742
$ws.Range("I136").Value = 2597.1304
$ws.Range("J136").Value = 3037.375
$ws.Range("K136").Value = 7791.3912
$ws.Range("L136").Value = 9112.125
$ws.Range("M136").Value = -5241.3912
$ws.Range("N136").Value = -14212.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1060.0197
$ws.Range("I58").Value = 1009.8461
$ws.Range("J58").Value = 1223.0834
$ws.Range("K58").Value = 1009.8461
$ws.Range("L58").Value = 1223.0834
$ws.Range("M58").Value = -806.8461
$ws.Range("N58").Value = -1629.0834
$ws.Range("H107").Value = 441.06897
$ws.Range("I107").Value = 421.3913
$ws.Range("J107").Value = 516.5
$ws.Range("K107").Value = 421.3913
$ws.Range("L107").Value = 516.5
$ws.Range("M107").Value = 1498.6087
$ws.Range("N107").Value = -4356.5
$ws.Range("H136").Value = 1060.0197
$ws.Range("I136").Value = 1009.8461
$ws.Range("J136").Value = 1223.0834
$ws.Range("K136").Value = 3029.5383
$ws.Range("L136").Value = 3669.2502
$ws.Range("M136").Value = -479.5383000000002
$ws.Range("N136").Value = -8769.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 900
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1246
$ws.Range("H19").Value = 9666
$ws.Range("J19").Value = 9666
$ws.Range("L19").Value = 28998
$ws.Range("N19").Value = -29346
$ws.Range("H29").Value = 832
$ws.Range("J29").Value = 1020
$ws.Range("L29").Value = 3060
$ws.Range("N29").Value = -3614
$ws.Range("H55").Value = 4645.533
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4645.533
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 13936.599
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -14290.599
$ws.Range("H58").Value = 1785
$ws.Range("J58").Value = 1785
$ws.Range("L58").Value = 5355
$ws.Range("N58").Value = -5611
$ws.Range("H68").Value = 162019.58
$ws.Range("I68").Value = 244457.22
$ws.Range("K68").Value = 733371.66
$ws.Range("M68").Value = -732560.66
$ws.Range("H71").Value = 162019.58
$ws.Range("I71").Value = 244457.22
$ws.Range("K71").Value = 2200114.98
$ws.Range("M71").Value = -2196058.98
$ws.Range("H122").Value = 1352.8948
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 1389.1666
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 12502.4994
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -17402.4994
$ws.Range("H131").Value = 2717.9092
$ws.Range("J131").Value = 3719.5527
$ws.Range("L131").Value = 11158.6581
$ws.Range("N131").Value = -21238.6581

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 27503.75
$ws.Range("H92").Value = 25025.5
$ws.Range("J92").Value = 25025.5
$ws.Range("L92").Value = 25025.5
$ws.Range("N92").Value = -28769.5
$ws.Range("H95").Value = 50344
$ws.Range("J95").Value = 50344
$ws.Range("L95").Value = 50344
$ws.Range("N95").Value = -55836
$ws.Range("H132").Value = 2410.484
$ws.Range("I132").Value = 1830.0476
$ws.Range("J132").Value = 3629.4
$ws.Range("K132").Value = 5490.142800000001
$ws.Range("L132").Value = 10888.2
$ws.Range("M132").Value = -2960.142800000001
$ws.Range("N132").Value = -15948.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 370
$ws.Range("I55").Value = 337.5
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 337.5
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -164.5
$ws.Range("N55").Value = -846
$ws.Range("H101").Value = 500000
$ws.Range("J101").Value = 500000
$ws.Range("L101").Value = 500000
$ws.Range("N101").Value = -506490
$ws.Range("H122").Value = 14291094
$ws.Range("I122").Value = 4964
$ws.Range("J122").Value = 22227834
$ws.Range("K122").Value = 14892
$ws.Range("L122").Value = 66683502
$ws.Range("M122").Value = -12442
$ws.Range("N122").Value = -66688402

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9200
$ws.Range("J69").Value = 9200
$ws.Range("L69").Value = 9200
$ws.Range("N69").Value = -10698
$ws.Range("H72").Value = 9200
$ws.Range("J72").Value = 9200
$ws.Range("L72").Value = 27600
$ws.Range("N72").Value = -35088
$ws.Range("H103").Value = 40150.5
$ws.Range("J103").Value = 40150.5
$ws.Range("L103").Value = 40150.5
$ws.Range("N103").Value = -42494.5
$ws.Range("H122").Value = 1499.9166
$ws.Range("I122").Value = 1501
$ws.Range("J122").Value = 1499.375
$ws.Range("K122").Value = 4503
$ws.Range("L122").Value = 4498.125
$ws.Range("M122").Value = -2053
$ws.Range("N122").Value = -9398.125

Write-Host "Applied 224 cell updates to Asura_Profits workbook"
